$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: insert a new, empty paragraph (its paragraph mark carries
# sz/szCs = 26, i.e. 13pt) immediately before the "Düh problémák
# gyermekekben és erőszakos játékok:" Heading 1 paragraph
# (TOC bookmark _Toc128345434).
# ------------------------------------------------------------------
$headingBookmark1 = $d.Bookmarks("_Toc128345434")
$headingStart1 = $headingBookmark1.Start
$insertionPoint = $d.Range($headingStart1, $headingStart1)

$newParaPkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$insertionPoint.InsertXML($newParaPkg)

# ------------------------------------------------------------------
# Change 2: delete the stand-alone page-break paragraph that sits
# right before the "Akkumulátorok hibái:" Heading 1 paragraph
# (TOC bookmark _Toc128345436) — it only contains a
# <w:lastRenderedPageBreak/><w:br w:type="page"/> run.
# ------------------------------------------------------------------
$headingBookmark2 = $d.Bookmarks("_Toc128345436")
$headingStart2 = $headingBookmark2.Start
$headingRange2 = $d.Range($headingStart2, $headingStart2)
$headingPara2 = $headingRange2.Paragraphs(1)
$pageBreakPara = $headingPara2.Previous()
$pageBreakPara.Range.Delete()
